$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 23, shifting the old row 23 (and below) down to row 24
$ws.Rows.Item(23).Insert()

# Populate the new row 23 with the updated record
$ws.Cells.Item(23, 1).Value = 10
$ws.Cells.Item(23, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(23, 3).Value = "La Araucanía"
$ws.Cells.Item(23, 4).Value = 44747
$ws.Cells.Item(23, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(23, 5).Value = 9
$ws.Cells.Item(23, 6).Value = 100112042
$ws.Cells.Item(23, 7).Value = "Locoto"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 80
$ws.Cells.Item(23, 11).Value = 2500
$ws.Cells.Item(23, 12).Value = 2500
$ws.Cells.Item(23, 13).Value = 2500
$ws.Cells.Item(23, 14).Value = "$/kilo"
$ws.Cells.Item(23, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(23, 16).Value = 2500
$ws.Cells.Item(23, 17).Value = 1
$ws.Cells.Item(23, 18).Value = "Hortaliza"
